$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 44 (PENDULINE ... row) and insert a copy above row 45 (REPAION-N row).
# This preserves styles/formatting exactly like the existing rows.
$ws.Rows.Item(44).Copy()
$ws.Rows.Item(45).Insert()

Write-Output "done"
